$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi following Dr Hou advice:
# add "ECs" as a new sending/target cluster (recomputed specificity values),
# expanding the LR-pairs table from 4 data rows to 6 data rows (rows 2-7).

# Row 2: ECs -> FAPs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna5"
$ws.Cells.Item(2, 3).Value = "Ephb2"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.230855
$ws.Cells.Item(2, 8).Value = 0.692565
$ws.Cells.Item(2, 9).Value = 0.06377305075821572
$ws.Cells.Item(2, 10).Value = 0.06377305075821572
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 6.346253666666667
$ws.Cells.Item(2, 14).Value = 19.038761
$ws.Cells.Item(2, 15).Value = 0.9446330608455225
$ws.Cells.Item(2, 16).Value = 0.9446330608455226
$ws.Cells.Item(2, 17).Value = 1.465064390218333
$ws.Cells.Item(2, 18).Value = 13.185579511965
$ws.Cells.Item(2, 19).Value = 0.06024213213719019
$ws.Cells.Item(2, 20).Value = 0.06024213213719019

# Row 3: ECs -> sCs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna5"
$ws.Cells.Item(3, 3).Value = "Ephb2"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.230855
$ws.Cells.Item(3, 8).Value = 0.692565
$ws.Cells.Item(3, 9).Value = 0.06377305075821572
$ws.Cells.Item(3, 10).Value = 0.06377305075821572
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.3719673333333333
$ws.Cells.Item(3, 14).Value = 1.115902
$ws.Cells.Item(3, 15).Value = 0.05536693915447755
$ws.Cells.Item(3, 16).Value = 0.05536693915447755
$ws.Cells.Item(3, 17).Value = 0.08587051873666667
$ws.Cells.Item(3, 18).Value = 0.7728346686299999
$ws.Cells.Item(3, 19).Value = 0.003530918621025538
$ws.Cells.Item(3, 20).Value = 0.003530918621025538

# Row 4: FAPs -> FAPs
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Efna5"
$ws.Cells.Item(4, 3).Value = "Ephb2"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.900731333333333
$ws.Cells.Item(4, 8).Value = 8.702194
$ws.Cells.Item(4, 9).Value = 0.8013189515350044
$ws.Cells.Item(4, 10).Value = 0.8013189515350045
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 6.346253666666667
$ws.Cells.Item(4, 14).Value = 19.038761
$ws.Cells.Item(4, 15).Value = 0.9446330608455225
$ws.Cells.Item(4, 16).Value = 0.9446330608455226
$ws.Cells.Item(4, 17).Value = 18.40877686018155
$ws.Cells.Item(4, 18).Value = 165.678991741634
$ws.Cells.Item(4, 19).Value = 0.7569523739020361
$ws.Cells.Item(4, 20).Value = 0.7569523739020363

# Row 5: FAPs -> sCs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efna5"
$ws.Cells.Item(5, 3).Value = "Ephb2"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.900731333333333
$ws.Cells.Item(5, 8).Value = 8.702194
$ws.Cells.Item(5, 9).Value = 0.8013189515350044
$ws.Cells.Item(5, 10).Value = 0.8013189515350045
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.3719673333333333
$ws.Cells.Item(5, 14).Value = 1.115902
$ws.Cells.Item(5, 15).Value = 0.05536693915447755
$ws.Cells.Item(5, 16).Value = 0.05536693915447755
$ws.Cells.Item(5, 17).Value = 1.078977298776444
$ws.Cells.Item(5, 18).Value = 9.710795688988
$ws.Cells.Item(5, 19).Value = 0.04436657763296833
$ws.Cells.Item(5, 20).Value = 0.04436657763296834

# Row 6: sCs -> FAPs
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Efna5"
$ws.Cells.Item(6, 3).Value = "Ephb2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.4883596666666667
$ws.Cells.Item(6, 8).Value = 1.465079
$ws.Cells.Item(6, 9).Value = 0.1349079977067798
$ws.Cells.Item(6, 10).Value = 0.1349079977067798
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 6.346253666666667
$ws.Cells.Item(6, 14).Value = 19.038761
$ws.Cells.Item(6, 15).Value = 0.9446330608455225
$ws.Cells.Item(6, 16).Value = 0.9446330608455226
$ws.Cells.Item(6, 17).Value = 3.099254325235445
$ws.Cells.Item(6, 18).Value = 27.893288927119
$ws.Cells.Item(6, 19).Value = 0.1274385548062961
$ws.Cells.Item(6, 20).Value = 0.1274385548062962

# Row 7: sCs -> sCs
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Efna5"
$ws.Cells.Item(7, 3).Value = "Ephb2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.4883596666666667
$ws.Cells.Item(7, 8).Value = 1.465079
$ws.Cells.Item(7, 9).Value = 0.1349079977067798
$ws.Cells.Item(7, 10).Value = 0.1349079977067798
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.3719673333333333
$ws.Cells.Item(7, 14).Value = 1.115902
$ws.Cells.Item(7, 15).Value = 0.05536693915447755
$ws.Cells.Item(7, 16).Value = 0.05536693915447755
$ws.Cells.Item(7, 17).Value = 0.1816538429175556
$ws.Cells.Item(7, 18).Value = 1.634884586258
$ws.Cells.Item(7, 19).Value = 0.007469442900483674
$ws.Cells.Item(7, 20).Value = 0.007469442900483674
